$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Excel "Find What"/"Replace With" constants:
#   LookAt: xlWhole = 1, xlPart = 2
#   SearchOrder: xlByRows = 1
# The trial/condition/filename labels embed the stimulus distance (D51,
# D64, D80) and size (S20, S25, S30) codes as substrings, so a partial
# (xlPart) match-and-replace across the whole used range reproduces the
# regenerated-order values from the commit: D51->D55, D64->D69, D80->D86,
# S30->S31 (S20/S25 stay as-is). None of the new codes are substrings of
# any of the old codes (or vice versa), so the four replacements are safe
# to run back-to-back in any order without double-substitution.
$used.Replace("D51", "D55", 2, 1, $false, $false)
$used.Replace("D64", "D69", 2, 1, $false, $false)
$used.Replace("D80", "D86", 2, 1, $false, $false)
$used.Replace("S30", "S31", 2, 1, $false, $false)
